$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Files")

# Shift existing F1 value ("*MISSING-FILE*") into G1, and set the new F1 header
$ws.Range("G1").Value = $ws.Range("F1").Value2
$ws.Range("F1").Value = "RELATION:contentLocation"

# Add new value for row 2
$ws.Range("F2").Value = "Catalina Park"

# Remove the old row 3 (CATALOG_pics.xlsx / 1) entirely
$ws.Rows.Item(3).Delete()
